$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived data for the Ptprz1-L1cam sending/target cluster matrix.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the various expression / specificity metrics (see header row).
$data = @(
    @("ECs","Ptprz1","L1cam","ECs",2,0.6666666666666666,0.06163433333333333,0.184903,0.006690894379667537,0.006690894379667537,3,1,7.741029,23.223087,0.4930486933812723,0.4930486933812723,0.4771131617289999,4.294018455561,0.003298936731447178,0.003298936731447178),
    @("ECs","Ptprz1","L1cam","FAPs",2,0.6666666666666666,0.06163433333333333,0.184903,0.006690894379667537,0.006690894379667537,3,1,0.3045986666666667,0.913796,0.0194007766416684,0.0194007766416684,0.01877373575422222,0.168963621788,0.0001298085473929243,0.0001298085473929243),
    @("ECs","Ptprz1","L1cam","MuSCs",2,0.6666666666666666,0.06163433333333333,0.184903,0.006690894379667537,0.006690894379667537,3,1,7.654706000000001,22.964118,0.4875505299770593,0.4875505299770593,0.4717927011726666,4.246134310554,0.003262149100827436,0.003262149100827436),
    @("FAPs","Ptprz1","L1cam","ECs",3,1,0.1030763333333333,0.309229,0.01118975126488057,0.01118975126488057,3,1,7.741029,23.223087,0.4930486933812723,0.4930486933812723,0.7979168855469999,7.181251969922999,0.005517092240410805,0.005517092240410805),
    @("FAPs","Ptprz1","L1cam","FAPs",3,1,0.1030763333333333,0.309229,0.01118975126488057,0.01118975126488057,3,1,0.3045986666666667,0.913796,0.0194007766416684,0.0194007766416684,0.03139691369822222,0.282572223284,0.0002170898649657745,0.0002170898649657745),
    @("FAPs","Ptprz1","L1cam","MuSCs",3,1,0.1030763333333333,0.309229,0.01118975126488057,0.01118975126488057,3,1,7.654706000000001,22.964118,0.4875505299770593,0.4875505299770593,0.7890190272246667,7.101171245022,0.005455569159503994,0.005455569159503994),
    @("MuSCs","Ptprz1","L1cam","ECs",3,1,9.046962666666667,27.140888,0.982119354355452,0.9821193543554518,3,1,7.741029,23.223087,0.4930486933812723,0.4930486933812723,70.032800364584,630.295203281256,0.4842326644094143,0.4842326644094143),
    @("MuSCs","Ptprz1","L1cam","FAPs",3,1,9.046962666666667,27.140888,0.982119354355452,0.9821193543554518,3,1,0.3045986666666667,0.913796,0.0194007766416684,0.0194007766416684,2.755692765649778,24.801234890848,0.0190538782293097,0.0190538782293097),
    @("MuSCs","Ptprz1","L1cam","MuSCs",3,1,9.046962666666667,27.140888,0.982119354355452,0.9821193543554518,3,1,7.654706000000001,22.964118,0.4875505299770593,0.4875505299770593,69.25183940630934,623.2665546567841,0.4788328117167279,0.4788328117167279)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $rowVals = $data[$i]
  $r = $startRow + $i
  for ($c = 1; $c -le $rowVals.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
  }
}
